$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1985.6
$ws.Range("I135").Value = 716.44446
$ws.Range("J135").Value = 3889.3333
$ws.Range("K135").Value = 6448.00014
$ws.Range("L135").Value = 35003.9997
$ws.Range("M135").Value = -3913.00014
$ws.Range("N135").Value = -40073.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7670.6206
$ws.Range("I32").Value = 7054.963
$ws.Range("K32").Value = 7054.963
$ws.Range("M32").Value = -6767.963
$ws.Range("H63").Value = 2396.875
$ws.Range("I63").Value = 1675
$ws.Range("J63").Value = 2830
$ws.Range("K63").Value = 1675
$ws.Range("L63").Value = 2830
$ws.Range("M63").Value = -989
$ws.Range("N63").Value = -4202
$ws.Range("H66").Value = 2396.875
$ws.Range("I66").Value = 1675
$ws.Range("J66").Value = 2830
$ws.Range("K66").Value = 8375
$ws.Range("L66").Value = 14150
$ws.Range("M66").Value = -4943
$ws.Range("N66").Value = -21014
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2089.0977
$ws.Range("I134").Value = 1745.2646
$ws.Range("J134").Value = 3759.1428
$ws.Range("K134").Value = 5235.793799999999
$ws.Range("L134").Value = 11277.4284
$ws.Range("M134").Value = -2700.793799999999
$ws.Range("N134").Value = -16347.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2658.842
$ws.Range("I58").Value = 1940.875
$ws.Range("J58").Value = 3181
$ws.Range("K58").Value = 1940.875
$ws.Range("L58").Value = 3181
$ws.Range("M58").Value = -1737.875
$ws.Range("N58").Value = -3587
$ws.Range("H132").Value = 5349.125
$ws.Range("I132").Value = 5297.9375
$ws.Range("K132").Value = 15893.8125
$ws.Range("M132").Value = -13363.8125
$ws.Range("H134").Value = 1539.4736
$ws.Range("I134").Value = 1486.3846
$ws.Range("J134").Value = 1654.5
$ws.Range("K134").Value = 4459.1538
$ws.Range("L134").Value = 4963.5
$ws.Range("M134").Value = -1924.1538
$ws.Range("N134").Value = -10033.5
$ws.Range("H136").Value = 2658.842
$ws.Range("I136").Value = 1940.875
$ws.Range("J136").Value = 3181
$ws.Range("K136").Value = 5822.625
$ws.Range("L136").Value = 9543
$ws.Range("M136").Value = -3272.625
$ws.Range("N136").Value = -14643
$ws.Range("H137").Value = 47500
$ws.Range("I137").Value = 30000
$ws.Range("K137").Value = 30000
$ws.Range("M137").Value = -24900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1193.9574
$ws.Range("I5").Value = 1084.1875
$ws.Range("K5").Value = 3252.5625
$ws.Range("M5").Value = -3140.5625
$ws.Range("H38").Value = 62.636364
$ws.Range("I38").Value = 43.333332
$ws.Range("J38").Value = 69.875
$ws.Range("K38").Value = 129.999996
$ws.Range("L38").Value = 209.625
$ws.Range("M38").Value = 217.000004
$ws.Range("N38").Value = -903.625
$ws.Range("H97").Value = 1360.7
$ws.Range("I97").Value = 1148.8
$ws.Range("K97").Value = 3446.4
$ws.Range("M97").Value = -2950.4
$ws.Range("H131").Value = 768.1900000000001
$ws.Range("I131").Value = 299.8889
$ws.Range("J131").Value = 814.5055
$ws.Range("K131").Value = 899.6667
$ws.Range("L131").Value = 2443.5165
$ws.Range("M131").Value = 4140.3333
$ws.Range("N131").Value = -12523.5165
$ws.Range("H132").Value = 1689.2572
$ws.Range("I132").Value = 761.4
$ws.Range("J132").Value = 2060.4
$ws.Range("K132").Value = 6852.599999999999
$ws.Range("L132").Value = 18543.6
$ws.Range("M132").Value = -4322.599999999999
$ws.Range("N132").Value = -23603.6
$ws.Range("H133").Value = 3062.6667
$ws.Range("I133").Value = 1806.6666
$ws.Range("J133").Value = 4946.6665
$ws.Range("K133").Value = 5419.9998
$ws.Range("L133").Value = 14839.9995
$ws.Range("M133").Value = -359.9997999999996
$ws.Range("N133").Value = -24959.9995
$ws.Range("H135").Value = 1193.9574
$ws.Range("I135").Value = 1084.1875
$ws.Range("K135").Value = 9757.6875
$ws.Range("M135").Value = -7222.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 43421.5
$ws.Range("I70").Value = 57303.844
$ws.Range("J70").Value = 5740.857
$ws.Range("K70").Value = 57303.844
$ws.Range("L70").Value = 5740.857
$ws.Range("M70").Value = -57033.844
$ws.Range("N70").Value = -6280.857
$ws.Range("H73").Value = 43421.5
$ws.Range("I73").Value = 57303.844
$ws.Range("J73").Value = 5740.857
$ws.Range("K73").Value = 57303.844
$ws.Range("L73").Value = 5740.857
$ws.Range("M73").Value = -56367.844
$ws.Range("N73").Value = -7612.857
$ws.Range("H122").Value = 1349.5
$ws.Range("I122").Value = 998.6667
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 2996.0001
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -546.0001000000002
$ws.Range("N122").Value = -9580
$ws.Range("H132").Value = 2992.4883
$ws.Range("I132").Value = 2562.3225
$ws.Range("J132").Value = 4103.75
$ws.Range("K132").Value = 7686.967500000001
$ws.Range("L132").Value = 12311.25
$ws.Range("M132").Value = -5156.967500000001
$ws.Range("N132").Value = -17371.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3503.7097
$ws.Range("I7").Value = 3844.375
$ws.Range("J7").Value = 3140.3333
$ws.Range("K7").Value = 3844.375
$ws.Range("L7").Value = 3140.3333
$ws.Range("M7").Value = -3732.375
$ws.Range("N7").Value = -3364.3333
$ws.Range("H22").Value = 2270.3125
$ws.Range("I22").Value = 2720
$ws.Range("K22").Value = 2720
$ws.Range("M22").Value = -2425
$ws.Range("H27").Value = 2270.3125
$ws.Range("I27").Value = 2720
$ws.Range("K27").Value = 2720
$ws.Range("M27").Value = -2613
$ws.Range("H46").Value = 1687483.1
$ws.Range("I46").Value = 490
$ws.Range("J46").Value = 2024881.8
$ws.Range("K46").Value = 490
$ws.Range("L46").Value = 2024881.8
$ws.Range("M46").Value = -302
$ws.Range("N46").Value = -2025257.8
$ws.Range("H126").Value = 3503.7097
$ws.Range("I126").Value = 3844.375
$ws.Range("J126").Value = 3140.3333
$ws.Range("K126").Value = 11533.125
$ws.Range("L126").Value = 9420.999899999999
$ws.Range("M126").Value = -9063.125
$ws.Range("N126").Value = -14360.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1638.9
$ws.Range("I126").Value = 1376.5555
$ws.Range("K126").Value = 4129.666499999999
$ws.Range("M126").Value = -1659.666499999999
